$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column changes: row 10 and row 33 change from "Data 2019" (53) to "NESCENT" (54)
$ws.Range("D10").Value = "NESCENT"
$ws.Range("D33").Value = "NESCENT"

# Fill in p-value columns E-J for rows with newly added data
# Row 10
$ws.Range("E10").Value = 0.624
$ws.Range("F10").Value = 0.272
$ws.Range("G10").Value = "<0.01"
$ws.Range("H10").Value = 0.0176
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = "<0.01"

# Row 12
$ws.Range("E12").Value = "<0.01"
$ws.Range("F12").Value = 0.0664
$ws.Range("G12").Value = "<0.01"
$ws.Range("H12").Value = "<0.01"
$ws.Range("I12").Value = 0.35
$ws.Range("J12").Value = "<0.01"

# Row 13
$ws.Range("E13").Value = 0.777
$ws.Range("F13").Value = 0.298
$ws.Range("G13").Value = "<0.01"
$ws.Range("H13").Value = "<0.01"
$ws.Range("I13").Value = 0.129
$ws.Range("J13").Value = "<0.01"

# Row 21
$ws.Range("E21").Value = 0.0783
$ws.Range("F21").Value = 0.477
$ws.Range("G21").Value = "<0.01"
$ws.Range("H21").Value = 0.0847
$ws.Range("I21").Value = 0.545
$ws.Range("J21").Value = "<0.01"

# Row 22
$ws.Range("E22").Value = 0.0256
$ws.Range("F22").Value = 0.625
$ws.Range("G22").Value = "<0.01"
$ws.Range("H22").Value = 0.407
$ws.Range("I22").Value = 0.204
$ws.Range("J22").Value = "<0.01"

# Row 23
$ws.Range("E23").Value = 0.0336
$ws.Range("F23").Value = 0.567
$ws.Range("G23").Value = "<0.01"
$ws.Range("H23").Value = 0.02
$ws.Range("I23").Value = 0.177
$ws.Range("J23").Value = "<0.01"

# Row 24
$ws.Range("E24").Value = "<0.01"
$ws.Range("F24").Value = 0.288
$ws.Range("G24").Value = "<0.01"
$ws.Range("H24").Value = 0.0244
$ws.Range("I24").Value = 0.126
$ws.Range("J24").Value = "<0.01"

# Row 25
$ws.Range("E25").Value = "<0.01"
$ws.Range("F25").Value = 0.926
$ws.Range("G25").Value = "<0.01"
$ws.Range("H25").Value = "<0.01"
$ws.Range("I25").Value = 0.488
$ws.Range("J25").Value = "<0.01"

# Row 26
$ws.Range("E26").Value = 0.0158
$ws.Range("F26").Value = 0.12
$ws.Range("G26").Value = "<0.01"
$ws.Range("H26").Value = "<0.01"
$ws.Range("I26").Value = 0.897
$ws.Range("J26").Value = "<0.01"

# Row 27
$ws.Range("E27").Value = 0.0115
$ws.Range("F27").Value = "<0.01"
$ws.Range("G27").Value = 0.881
$ws.Range("H27").Value = 0.0363
$ws.Range("I27").Value = "<0.01"
$ws.Range("J27").Value = 1

# Row 28
$ws.Range("E28").Value = "<0.01"
$ws.Range("F28").Value = 0.056
$ws.Range("G28").Value = "<0.01"
$ws.Range("H28").Value = 0.143
$ws.Range("I28").Value = "<0.01"
$ws.Range("J28").Value = "<0.01"

# Row 29
$ws.Range("E29").Value = "<0.01"
$ws.Range("F29").Value = 0.193
$ws.Range("G29").Value = "<0.01"
$ws.Range("H29").Value = "<0.01"
$ws.Range("I29").Value = 0.0174
$ws.Range("J29").Value = "<0.01"

# Row 30
$ws.Range("E30").Value = "<0.01"
$ws.Range("F30").Value = "<0.01"
$ws.Range("G30").Value = "<0.01"
$ws.Range("H30").Value = "<0.01"
$ws.Range("I30").Value = 0.0426
$ws.Range("J30").Value = "<0.01"

# Row 31
$ws.Range("E31").Value = "<0.01"
$ws.Range("F31").Value = "<0.01"
$ws.Range("G31").Value = 0.0852
$ws.Range("H31").Value = 0.0695
$ws.Range("I31").Value = 0.301
$ws.Range("J31").Value = 0.237

# Row 32
$ws.Range("E32").Value = 0.0103
$ws.Range("F32").Value = "<0.01"
$ws.Range("G32").Value = "<0.01"
$ws.Range("H32").Value = 0.0211
$ws.Range("I32").Value = 0.0678
$ws.Range("J32").Value = "<0.01"

# Row 33
$ws.Range("E33").Value = "<0.01"
$ws.Range("F33").Value = 0.43
$ws.Range("G33").Value = "<0.01"
$ws.Range("H33").Value = "<0.01"
$ws.Range("I33").Value = 0.939
$ws.Range("J33").Value = "<0.01"

# Row 34
$ws.Range("E34").Value = 0.195
$ws.Range("F34").Value = 0.0873
$ws.Range("G34").Value = "<0.01"
$ws.Range("H34").Value = 0.258
$ws.Range("I34").Value = 0.141
$ws.Range("J34").Value = "<0.01"

# Row 35
$ws.Range("E35").Value = 0.0335
$ws.Range("F35").Value = 0.93
$ws.Range("G35").Value = "<0.01"
$ws.Range("H35").Value = 0.0433
$ws.Range("I35").Value = "<0.01"
$ws.Range("J35").Value = "<0.01"

# Row 36
$ws.Range("E36").Value = 1
$ws.Range("F36").Value = 0.908
$ws.Range("G36").Value = 0.667
$ws.Range("H36").Value = 0.0429
$ws.Range("I36").Value = 0.316
$ws.Range("J36").Value = 0.175

# Row 37
$ws.Range("E37").Value = 0.232
$ws.Range("F37").Value = 0.276
$ws.Range("G37").Value = 0.85
$ws.Range("H37").Value = 0.369
$ws.Range("I37").Value = 0.0277
$ws.Range("J37").Value = 1

# Row 38
$ws.Range("E38").Value = 0.593
$ws.Range("F38").Value = 0.798
$ws.Range("G38").Value = "<0.01"
$ws.Range("H38").Value = 0.492
$ws.Range("I38").Value = 0.692
$ws.Range("J38").Value = "<0.01"

# Row 39
$ws.Range("E39").Value = 0.616
$ws.Range("F39").Value = 0.654
$ws.Range("G39").Value = 0.441
$ws.Range("H39").Value = 0.412
$ws.Range("I39").Value = 0.0733
$ws.Range("J39").Value = 0.838

# Row 40
$ws.Range("E40").Value = 0.226
$ws.Range("F40").Value = 0.933
$ws.Range("G40").Value = 0.249
$ws.Range("H40").Value = 0.0575
$ws.Range("I40").Value = 0.727
$ws.Range("J40").Value = 0.466

# Row 41
$ws.Range("E41").Value = 0.398
$ws.Range("F41").Value = 0.602
$ws.Range("G41").Value = 0.0905
$ws.Range("H41").Value = 0.609
$ws.Range("I41").Value = 0.857
$ws.Range("J41").Value = 0.544

# Row 42
$ws.Range("E42").Value = 0.361
$ws.Range("F42").Value = 0.0321
$ws.Range("G42").Value = 0.826
$ws.Range("H42").Value = 0.21
$ws.Range("I42").Value = 0.0483
$ws.Range("J42").Value = 0.742

# Row 43
$ws.Range("E43").Value = "<0.01"
$ws.Range("F43").Value = 0.763
$ws.Range("G43").Value = "<0.01"
$ws.Range("H43").Value = "<0.01"
$ws.Range("I43").Value = 0.308
$ws.Range("J43").Value = "<0.01"

# Row 44
$ws.Range("E44").Value = "<0.01"
$ws.Range("F44").Value = 0.755
$ws.Range("G44").Value = "<0.01"
$ws.Range("H44").Value = "<0.01"
$ws.Range("I44").Value = 0.801
$ws.Range("J44").Value = "<0.01"

# Update selection to J37 (also clears the scrolled topLeftCell)
$ws.Range("J37").Select()
